$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 38571.285
$ws.Range("I33").Value = 56423.223
$ws.Range("J33").Value = 6437.8
$ws.Range("K33").Value = 56423.223
$ws.Range("L33").Value = 6437.8
$ws.Range("M33").Value = -56194.223
$ws.Range("N33").Value = -6895.8

$ws.Range("H38").Value = 2562.4285
$ws.Range("J38").Value = 5696.778
$ws.Range("L38").Value = 17090.334
$ws.Range("N38").Value = -17834.334

$ws.Range("H58").Value = 700
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H64").Value = 4496.4863
$ws.Range("I64").Value = 4738.4
$ws.Range("K64").Value = 4738.4
$ws.Range("M64").Value = -4490.4

$ws.Range("H67").Value = 4496.4863
$ws.Range("I67").Value = 4738.4
$ws.Range("K67").Value = 4738.4
$ws.Range("M67").Value = -3880.4

$ws.Range("H87").Value = 84345.92999999999
$ws.Range("J87").Value = 86082
$ws.Range("L87").Value = 86082
$ws.Range("N87").Value = -88578

$ws.Range("H90").Value = 84345.92999999999
$ws.Range("J90").Value = 86082
$ws.Range("L90").Value = 258246
$ws.Range("N90").Value = -270726

$ws.Range("H113").Value = 3605.125
$ws.Range("J113").Value = 4158.3
$ws.Range("L113").Value = 4158.3
$ws.Range("N113").Value = -10666.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9191.322
$ws.Range("I32").Value = 3390.6072
$ws.Range("K32").Value = 3390.6072
$ws.Range("M32").Value = -3103.6072

$ws.Range("H45").Value = 1948.5333
$ws.Range("J45").Value = 3249.5
$ws.Range("L45").Value = 3249.5
$ws.Range("N45").Value = -4003.5

$ws.Range("H63").Value = 4987.2856
$ws.Range("J63").Value = 9332.666999999999
$ws.Range("L63").Value = 9332.666999999999
$ws.Range("N63").Value = -10704.667

$ws.Range("H66").Value = 4987.2856
$ws.Range("J66").Value = 9332.666999999999
$ws.Range("L66").Value = 46663.335
$ws.Range("N66").Value = -53527.335

$ws.Range("H132").Value = 2588.4285
$ws.Range("I132").Value = 1403.92
$ws.Range("J132").Value = 5549.7
$ws.Range("K132").Value = 4211.76
$ws.Range("L132").Value = 16649.1
$ws.Range("M132").Value = -1681.76
$ws.Range("N132").Value = -21709.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 449.57144
$ws.Range("I22").Value = 453.14285
$ws.Range("K22").Value = 453.14285
$ws.Range("M22").Value = -280.14285

$ws.Range("H75").Value = 155270
$ws.Range("J75").Value = 179125
$ws.Range("L75").Value = 179125
$ws.Range("N75").Value = -180997

$ws.Range("H78").Value = 155270
$ws.Range("J78").Value = 179125
$ws.Range("L78").Value = 537375
$ws.Range("N78").Value = -546735

$ws.Range("H82").Value = 152222.22
$ws.Range("I82").Value = 30257
$ws.Range("K82").Value = 30257
$ws.Range("M82").Value = -29874

$ws.Range("H85").Value = 152222.22
$ws.Range("I85").Value = 30257
$ws.Range("K85").Value = 30257
$ws.Range("M85").Value = -28931

$ws.Range("H86").Value = 1265.35
$ws.Range("I86").Value = 764
$ws.Range("J86").Value = 1878.1111
$ws.Range("K86").Value = 764
$ws.Range("L86").Value = 1878.1111
$ws.Range("M86").Value = 359
$ws.Range("N86").Value = -4124.1111

$ws.Range("H89").Value = 1265.35
$ws.Range("I89").Value = 764
$ws.Range("J89").Value = 1878.1111
$ws.Range("K89").Value = 3820
$ws.Range("L89").Value = 9390.5555
$ws.Range("M89").Value = 1796
$ws.Range("N89").Value = -20622.5555

$ws.Range("H105").Value = 2242.7083
$ws.Range("I105").Value = 2065.9546
$ws.Range("J105").Value = 4187
$ws.Range("K105").Value = 2065.9546
$ws.Range("L105").Value = 4187
$ws.Range("M105").Value = -318.9546
$ws.Range("N105").Value = -7681

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 98956.39999999999
$ws.Range("I69").Value = 35880
$ws.Range("K69").Value = 35880
$ws.Range("M69").Value = -35131

$ws.Range("H72").Value = 98956.39999999999
$ws.Range("I72").Value = 35880
$ws.Range("K72").Value = 107640
$ws.Range("M72").Value = -103896

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 55
$ws.Range("I6").Value = 72.5
$ws.Range("J6").Value = 20
$ws.Range("K6").Value = 217.5
$ws.Range("L6").Value = 60
$ws.Range("M6").Value = -104.5
$ws.Range("N6").Value = -286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 146247
$ws.Range("J106").Value = 146247
$ws.Range("L106").Value = 146247
$ws.Range("N106").Value = -148771

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 89300.086
$ws.Range("I7").Value = 105610.1
$ws.Range("J7").Value = 7750
$ws.Range("K7").Value = 105610.1
$ws.Range("L7").Value = 7750
$ws.Range("M7").Value = -105498.1
$ws.Range("N7").Value = -7974

$ws.Range("H64").Value = 173291.5
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 173291.5
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 173291.5
$ws.Range("N64").Value = -173741.5
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 173291.5
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 173291.5
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 173291.5
$ws.Range("N67").Value = -174851.5
$ws.Range("M67").ClearContents()

$ws.Range("H93").Value = 1541.5
$ws.Range("I93").Value = 1300.875
$ws.Range("J93").Value = 1862.3334
$ws.Range("K93").Value = 1300.875
$ws.Range("L93").Value = 1862.3334
$ws.Range("M93").Value = -52.875
$ws.Range("N93").Value = -4358.3334

$ws.Range("H126").Value = 89300.086
$ws.Range("I126").Value = 105610.1
$ws.Range("J126").Value = 7750
$ws.Range("K126").Value = 316830.3
$ws.Range("L126").Value = 23250
$ws.Range("M126").Value = -314360.3
$ws.Range("N126").Value = -28190

$ws.Range("H132").Value = 7754.9
$ws.Range("I132").Value = 7923.636
$ws.Range("J132").Value = 7548.6665
$ws.Range("K132").Value = 23770.908
$ws.Range("L132").Value = 22645.9995
$ws.Range("M132").Value = -21240.908
$ws.Range("N132").Value = -27705.9995
